# "add property for npc"
# The NPC.xlsx "Property" sheet defines one NPC property per row (columns
# A:J = Id, Type, Public, Private, Save, View, Index, SaveInterval,
# RelationValue, Desc). This adds a new "Height" (模型高度) property as the
# next row, right after the existing "SkillIDRef" (技能列表索引) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

$newRow = $ws.UsedRange.Rows.Count + 1

$ws.Cells.Item($newRow, 1).Value = "Height"
$ws.Cells.Item($newRow, 2).Value = "float"
$ws.Cells.Item($newRow, 3).Value = $false
$ws.Cells.Item($newRow, 4).Value = $false
$ws.Cells.Item($newRow, 5).Value = $false
$ws.Cells.Item($newRow, 6).Value = $true
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = "Friend"
$ws.Cells.Item($newRow, 9).NumberFormat = "@"
$ws.Cells.Item($newRow, 10).Value = "模型高度"

$null = $ws.Cells.Item($newRow, 10).Select()
